$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation for 2026/02/12 (木, hour 17, ranking 201) was recorded.
# It belongs right after the existing 2026/02/12 rows (805-810) and before
# the 2026/12/29 block, so insert a fresh row at 811 and push everything
# that was there (811..852) down to (812..853).
$ws.Rows.Item(811).Insert()

# Column A holds dates formatted/stored as plain text (e.g. "2026/12/29"),
# so force Text formatting before assigning to stop Excel from silently
# reinterpreting the string as a date serial, then drop back to the
# workbook's normal (unstyled) look so the new row matches its neighbors.
$ws.Cells.Item(811, 1).NumberFormat = "@"
$ws.Cells.Item(811, 1).Value = "2026/02/12"
$ws.Cells.Item(811, 1).Style = "Normal"

$ws.Cells.Item(811, 2).Value = "木"
$ws.Cells.Item(811, 3).Value = 17
$ws.Cells.Item(811, 4).Value = 201
